$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - "Starting with Coin System/Remove Coin Function"
$ws.Range("D2").Value = "6h"
$ws.Range("E2").Value = "-"

# Row 3 - "daily login and shop view"
$ws.Range("D3").Value = "4h"
$ws.Range("E3").Value = "-"

# Row 6 - "Item Storage"
$ws.Range("D6").Value = "3h"
$ws.Range("E6").Value = "-"
